$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.22"
$ws.Range("E2").Value = "'1.59%"
$ws.Range("D3").Value = "'37.33"
$ws.Range("E3").Value = "'0.03%"
$ws.Range("D4").Value = "'5.150"
$ws.Range("E4").Value = "'1.29%"
$ws.Range("D5").Value = "'0.07813"
$ws.Range("E5").Value = "'1.23%"
$ws.Range("D6").Value = "'4.423"
$ws.Range("D7").Value = "'1.905"
$ws.Range("E7").Value = "'0.96%"
$ws.Range("D8").Value = "'8.269"
$ws.Range("E8").Value = "'0.84%"
$ws.Range("D9").Value = "'2.794"
$ws.Range("E9").Value = "'-7.41%"
$ws.Range("D10").Value = "'0.9198"
$ws.Range("E10").Value = "'0.22%"
$ws.Range("D11").Value = "'0.1189"
$ws.Range("E11").Value = "'3.42%"
$ws.Range("D12").Value = "'0.1923"
$ws.Range("E12").Value = "'2.24%"
$ws.Range("D13").Value = "'0.09045"
$ws.Range("E13").Value = "'3.92%"
$ws.Range("D14").Value = "'0.03354"
$ws.Range("E14").Value = "'-1.40%"
$ws.Range("D15").Value = "'0.09611"
$ws.Range("E15").Value = "'-0.95%"
$ws.Range("D16").Value = "'0.001399"
$ws.Range("E16").Value = "'1.72%"
$ws.Range("D17").Value = "'0.005781"
$ws.Range("E17").Value = "'-2.20%"
$ws.Range("D18").Value = "'3.519"
$ws.Range("E18").Value = "'-1.86%"
$ws.Range("D19").Value = "'0.3441"
$ws.Range("E19").Value = "'0.98%"
$ws.Range("D20").Value = "'5.267"
$ws.Range("E20").Value = "'5.00%"
$ws.Range("D21").Value = "'0.1272"
$ws.Range("E21").Value = "'-1.32%"
$ws.Range("D22").Value = "'0.2596"
$ws.Range("E22").Value = "'0.01%"
$ws.Range("D23").Value = "'0.04382"
$ws.Range("E23").Value = "'1.46%"
$ws.Range("D24").Value = "'0.001255"
$ws.Range("E24").Value = "'3.36%"
$ws.Range("D25").Value = "'0.004663"
$ws.Range("E25").Value = "'2.81%"
$ws.Range("D26").Value = "'0.0001366"
$ws.Range("E26").Value = "'1.08%"
$ws.Range("D27").Value = "'0.0004008"
$ws.Range("D39").Value = "'0.02286"
$ws.Range("E39").Value = "'3.73%"
$ws.Range("D40").Value = "'0.05052"
$ws.Range("E40").Value = "'2.72%"
$ws.Range("D41").Value = "'0.007472"
$ws.Range("E41").Value = "'-1.18%"
$ws.Range("D42").Value = "'0.009095"
$ws.Range("E42").Value = "'-8.48%"
$ws.Range("D43").Value = "'0.1350"
$ws.Range("E43").Value = "'1.26%"
$ws.Range("D44").Value = "'0.001959"
$ws.Range("E44").Value = "'-1.84%"
$ws.Range("D45").Value = "'0.009435"
$ws.Range("E45").Value = "'11.14%"
$ws.Range("D46").Value = "'0.00006701"
$ws.Range("E46").Value = "'2.35%"
$ws.Range("E47").Value = "'0.31%"
$ws.Range("D48").Value = "'0.003299"
$ws.Range("E48").Value = "'10.01%"
$ws.Range("D49").Value = "'0.001005"
$ws.Range("E49").Value = "'-22.85%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.31%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.31%"
